# Update "想去人数" (interest-count) figures in column F on the
# 展览 (Exhibition) and 全部类型 (All types) sheets to the latest scrape.
#
# 展览 sheet rows -> 全部类型 sheet rows (same events, different row offsets)
#   F6  : 206   -> 208
#   F7  : 780   -> 781
#   F14 : 6602  -> 6605   (全部类型 F16)
#   F20 : 1043  -> 1047   (全部类型 F23)
#   F21 : 15815 -> 15818  (全部类型 F24)
#   F27 : 11192 -> 11194  (全部类型 F31)
#   F29 : 4389  -> 4390   (全部类型 F33)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value map for each sheet.
$exhibitionUpdates = @{
    6  = 208
    7  = 781
    14 = 6605
    20 = 1047
    21 = 15818
    27 = 11194
    29 = 4390
}

$allTypesUpdates = @{
    6  = 208
    7  = 781
    16 = 6605
    23 = 1047
    24 = 15818
    31 = 11194
    33 = 4390
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
